$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Debtor info block
$ws.Range("B11").Value = "Pablo"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "1234567890"

# Concept
$ws.Range("B23").Value = "Concepto de prueba"

# Question b: "Soy declarante del Impuesto de Renta" -> mark SI (F41) instead of NO (F42)
$ws.Range("F42").Value = ""
$ws.Range("F41").Value = "x"

# Question c: "Soy residente fiscal en Colombia" -> mark SI (F44) instead of NO (F45)
$ws.Range("F45").Value = ""
$ws.Range("F44").Value = "x"

# City and date
$ws.Range("C48").Value = "Bogota"
$ws.Range("E48").Value = (Get-Date -Year 2023 -Month 4 -Day 1 -Hour 0 -Minute 0 -Second 0)

# Signature block
$ws.Range("C51").Value = "Pablo"
$ws.Range("C52").NumberFormat = "@"
$ws.Range("C52").Value = "1234567890"
$ws.Range("C53").Value = "Calle 123"
$ws.Range("C54").NumberFormat = "@"
$ws.Range("C54").Value = "1234567890"

# Bank
$ws.Range("B58").Value = "Banco de Prueba"
